# Apply the StructureDefinition-local-race-cd.xlsx update:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to 2022-01-21T20:46:54+00:00
#  - set Publisher to "Alvearie Team"
#  - replace the duplicated "Contact / No display for ContactDetail" row with
#    "Jurisdiction / United States of America", dropping the extra duplicate row
#  - fill in the root Extension row's Short/Definition on the Elements sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Metadata" sheet
$ws2 = $wb.Worksheets.Item(2)   # "Elements" sheet

# --- Metadata sheet updates ---

# Version value (row 3)
$ws1.Range("B3").Value = "6.0.0"

# Date value (row 8)
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (row 9), previously blank
$ws1.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was an exact duplicate of the old row 10 ("Contact" / "No display for ContactDetail").
# Remove it entirely so everything below shifts up by one row.
$ws1.Rows.Item(11).Delete()

# --- Elements sheet updates ---
# Root Extension element row (row 2): Short / Definition columns (K / L)
$ws2.Range("K2").Value = "Local Race Code"
$ws2.Range("L2").Value = "Customer-specific code for the race of the person"
